$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 72 (pushes existing rows 72:87 down to 73:88),
# inheriting formatting from the row above (matches the existing date-style
# cell in column D).
$ws.Rows.Item(72).Insert()

# Populate the newly inserted row with the new weekly price-record data point.
$ws.Range("A72").Value = 5
$ws.Range("B72").Value = "Macroferia Regional de Talca"
$ws.Range("C72").Value = "Maule"
$ws.Range("D72").Value = 44504
$ws.Range("E72").Value = 7
$ws.Range("F72").Value = "Fruta"
$ws.Range("G72").Value = 100108
$ws.Range("H72").Value = "Tropicales y subtropicales"
$ws.Range("I72").Value = 100108002
$ws.Range("J72").Value = "Mango"
$ws.Range("K72").Value = "Sin especificar"
$ws.Range("L72").Value = "Primera"
$ws.Range("M72").Value = 260
$ws.Range("N72").Value = 7000
$ws.Range("O72").Value = 7000
$ws.Range("P72").Value = 7000
$ws.Range("Q72").Value = "$/bandeja 4 kilos"
$ws.Range("R72").Value = "Perú"
$ws.Range("S72").Value = 1750
$ws.Range("T72").Value = 4
